$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (46075 -> 46076) for every data row (rows 2 through 537).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 537 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 46076
